# Mod 6 last changes
# Slide 11 ("SASS" / "LESS" comparison slide):
#  1. Merge "Supported by Compass, " + "a " runs into a single run
#     "Supported by Compass, a ", and re-seat "hugely " as its own
#     (now dirty="0") run ahead of the untouched "popular " run.
#  2. Split "Conservative but wise about features" into
#     "Conservative (but wise) " + "about features".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Part 1: "Supported by Compass, a hugely popular mixin library" ---
$shCompass = $s.Shapes.Item(1)
$paraCompass = $shCompass.TextFrame.TextRange.Paragraphs(3)

# Remove the three old runs covering "Supported by Compass, a hugely "
# (chars 1-31) outright, then re-insert the merged wording right before
# the untouched "popular" run so the new text picks up a fresh,
# dirty="0" rPr instead of inheriting the old un-dirtied "a "/"hugely "
# run formatting.
$oldSpan = $paraCompass.Characters(1, 31)
$oldSpan.Delete() | Out-Null

$anchor = $paraCompass.Characters(1, 1)
$anchor.InsertBefore("Supported by Compass, a hugely ") | Out-Null

# Re-split the freshly inserted text at the "Supported by Compass, a "
# boundary (24 chars) so it forms its own run distinct from "hugely ".
$rCompass1 = $paraCompass.Characters(1, 24)
$rCompass1.Text = "Supported by Compass, a "

# The insertion point anchored on the old "p" of "popular", so that
# single character rode along with "hugely " (-> "hugely p") and with
# "opular " left behind as its own (untouched-formatting) run. Move it
# back so the runs read "hugely " / "popular " again.
$rCompass2 = $paraCompass.Characters(25, 7)
$rCompass2.Text = "hugely "

$rCompass3 = $paraCompass.Characters(32, 8)
$rCompass3.Text = "popular "

# --- Part 2: "Conservative but wise about features" ---
$shLess = $s.Shapes.Item(3)
$paraLess = $shLess.TextFrame.TextRange.Paragraphs(2)

# Rewrite the whole (single-run) sentence with the new wording first.
$fullLess = $paraLess.Characters(1, $paraLess.Length)
$fullLess.Text = "Conservative (but wise) about features"

# Now split it into two runs at the word boundary.
$rLess1 = $paraLess.Characters(1, 24)
$rLess1.Text = "Conservative (but wise) "

$rLess2 = $paraLess.Characters(25, 14)
$rLess2.Text = "about features"
